$d = $word.ActiveDocument

# "We will make use of API calls to enable our AI coach (more on that in
# Features). We will also make use of API calls to display the information..."
# -> "We will make use of API calls to display the information..."
#
# Removes the now-stale AI-coach reference (and the duplicated "We will
# also make use of API calls to" lead-in) from the API Calls paragraph,
# leaving a single, still-accurate sentence about API usage.
$found = $d.Content.Find.Execute(
    "enable our AI coach (more on that in Features). We will also make use of API calls to ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2
)
Write-Output "Replaced AI-coach reference: $found"
